$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6: date text in column A, numeric values in B:G
$ws.Range("A2").Value = "18/10/2023"
$ws.Range("B2").Value = 4512.4
$ws.Range("C2").Value = 4512.4
$ws.Range("D2").Value = 5012.4
$ws.Range("E2").Value = 5012.4
$ws.Range("F2").Value = 500
$ws.Range("G2").Value = 111.08
$ws.Range("A3").Value = "19/10/2023"
$ws.Range("B3").Value = 5475.3
$ws.Range("C3").Value = 9987.7
$ws.Range("D3").Value = 8451.2
$ws.Range("E3").Value = 13463.6
$ws.Range("F3").Value = 3475.9
$ws.Range("G3").Value = 134.8
$ws.Range("A4").Value = "20/10/2023"
$ws.Range("B4").Value = 8800
$ws.Range("C4").Value = 18787.7
$ws.Range("D4").Value = 8900
$ws.Range("E4").Value = 22363.6
$ws.Range("F4").Value = 3575.9
$ws.Range("G4").Value = 119.03
$ws.Range("A5").Value = "21/10/2023"
$ws.Range("B5").Value = 9000
$ws.Range("C5").Value = 27787.7
$ws.Range("D5").Value = 11000
$ws.Range("E5").Value = 33363.6
$ws.Range("F5").Value = 5575.9
$ws.Range("G5").Value = 120.07
$ws.Range("A6").Value = "22/10/2023"
$ws.Range("B6").Value = 9000.6
$ws.Range("C6").Value = 36788.3
$ws.Range("D6").Value = 15060.1
$ws.Range("E6").Value = 48423.7
$ws.Range("F6").Value = 11635.4
$ws.Range("G6").Value = 131.63
# Rows 7-8: all values stored as literal text (e.g. "4800.00"), so force
# text formatting before assignment to stop numeric-looking strings from
# being auto-converted to numbers, then restore the default "Normal"
# style so the cells keep no explicit style reference (just like the
# other new rows).
$ws.Range("A7:G7").NumberFormat = "@"
$ws.Range("A7").Value = "23/10/2023"
$ws.Range("B7").Value = "4800.00"
$ws.Range("C7").Value = "41588.30"
$ws.Range("D7").Value = "9220.00"
$ws.Range("E7").Value = "57643.70"
$ws.Range("F7").Value = "16055.40"
$ws.Range("G7").Value = "138.61"
$ws.Range("A7:G7").Style = "Normal"
$ws.Range("A8:G8").NumberFormat = "@"
$ws.Range("A8").Value = "25/10/2023"
$ws.Range("B8").Value = "7000.00"
$ws.Range("C8").Value = "48588.30"
$ws.Range("D8").Value = "4650.08"
$ws.Range("E8").Value = "62293.78"
$ws.Range("F8").Value = "13705.48"
$ws.Range("G8").Value = "128.21"
$ws.Range("A8:G8").Style = "Normal"
